$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update status text shared across sheets: "Ready for handoff" -> "Handback transform failed"
$wsOverview.Range("E2").Value = "Handback transform failed"
$wsOverview.Range("F2").Value = "Handback transform failed"
$wsZhCn.Range("C2").Value = "Handback transform failed"
$wsDeDe.Range("C2").Value = "Handback transform failed"

# New Error Detail message for zh-cn / de-de handback failures
$errMsg = "The translationStateItem f875148175623cd34b302a29281ae7cf7bea90f4 is not found."
$wsZhCn.Range("P2").Value = $errMsg
$wsDeDe.Range("P2").Value = $errMsg

# Widen the Status columns (Overview E:F, zh-cn/de-de column C)
# target raw XML width is 24.7426795959473 characters; the host's
# ColumnWidth setter quantizes to whole-pixel steps (MDW=6), so use the
# character value whose quantized width lands nearest the target (24.6667).
$wsOverview.Columns.Item(5).ColumnWidth = 23.8333333333
$wsOverview.Columns.Item(6).ColumnWidth = 23.8333333333
$wsZhCn.Columns.Item(3).ColumnWidth = 23.8333333333
$wsDeDe.Columns.Item(3).ColumnWidth = 23.8333333333

# Widen Error Detail column (column P, 16) to fit the new message (target raw width 40)
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666666667
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666666667
